# Updates the cryptos price/volume table (columns D and E, rows 2-51) with
# freshly scraped values. All D/E cells in this sheet are stored as plain
# text (many "prices" use dotted thousands-grouping like "30.365.71", and
# the "Volume(1h)" column is a padded "  +n.nn%  " string) so every write
# below forces a text number-format before assigning the value and then
# restores the cell's style to "Normal" afterwards. This keeps the cell's
# underlying type as text (matching the original inlineStr cells) instead
# of letting the host auto-coerce number-looking strings (e.g. "1.002")
# into numeric values, while leaving no residual per-cell style applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "30.365.71" },
    @{ Cell = "E2"; Value = "  +1.48%  " },
    @{ Cell = "D3"; Value = "2.011.15" },
    @{ Cell = "E3"; Value = "  +4.57%  " },
    @{ Cell = "D4"; Value = "1.002" },
    @{ Cell = "E4"; Value = "  +0.21%  " },
    @{ Cell = "D5"; Value = "324.91" },
    @{ Cell = "E5"; Value = "  +1.43%  " },
    @{ Cell = "E6"; Value = "  +0.13%  " },
    @{ Cell = "D7"; Value = "0.5134" },
    @{ Cell = "E7"; Value = "  +1.46%  " },
    @{ Cell = "D8"; Value = "0.4256" },
    @{ Cell = "E8"; Value = "  +5.61%  " },
    @{ Cell = "D9"; Value = "0.08744" },
    @{ Cell = "E9"; Value = "  +4.82%  " },
    @{ Cell = "D10"; Value = "43.56" },
    @{ Cell = "E10"; Value = "  +2.47%  " },
    @{ Cell = "E11"; Value = "  +2.67%  " },
    @{ Cell = "D12"; Value = "24.47" },
    @{ Cell = "E12"; Value = "  +2.91%  " },
    @{ Cell = "E13"; Value = "  +4.73%  " },
    @{ Cell = "D14"; Value = "6.634" },
    @{ Cell = "E14"; Value = "  +3.50%  " },
    @{ Cell = "D15"; Value = "7.452" },
    @{ Cell = "E15"; Value = "  +2.99%  " },
    @{ Cell = "D16"; Value = "1.002" },
    @{ Cell = "E16"; Value = "  +0.19%  " },
    @{ Cell = "E17"; Value = "  +2.16%  " },
    @{ Cell = "E18"; Value = "  +1.37%  " },
    @{ Cell = "D19"; Value = "0.06541" },
    @{ Cell = "E19"; Value = "  +0.47%  " },
    @{ Cell = "D20"; Value = "18.84" },
    @{ Cell = "E20"; Value = "  +3.05%  " },
    @{ Cell = "E21"; Value = "  +0.03%  " },
    @{ Cell = "D22"; Value = "6.204" },
    @{ Cell = "E22"; Value = "  +4.28%  " },
    @{ Cell = "D23"; Value = "30.431.85" },
    @{ Cell = "E23"; Value = "  +1.65%  " },
    @{ Cell = "D24"; Value = "11.84" },
    @{ Cell = "D25"; Value = "2.252" },
    @{ Cell = "E25"; Value = "  +2.74%  " },
    @{ Cell = "D26"; Value = "2.251.86" },
    @{ Cell = "E26"; Value = "  +5.16%  " },
    @{ Cell = "E27"; Value = "  +1.28%  " },
    @{ Cell = "D28"; Value = "161.92" },
    @{ Cell = "E28"; Value = "  -0.08%  " },
    @{ Cell = "D29"; Value = "2.437" },
    @{ Cell = "E29"; Value = "  +4.82%  " },
    @{ Cell = "D30"; Value = "131.29" },
    @{ Cell = "E30"; Value = "  +1.78%  " },
    @{ Cell = "D31"; Value = "1.144" },
    @{ Cell = "E31"; Value = "  +1.16%  " },
    @{ Cell = "E32"; Value = "  +1.71%  " },
    @{ Cell = "D33"; Value = "6.091" },
    @{ Cell = "E33"; Value = "  +2.10%  " },
    @{ Cell = "E34"; Value = "  +1.07%  " },
    @{ Cell = "D35"; Value = "1.360" },
    @{ Cell = "E35"; Value = "  +13.80%  " },
    @{ Cell = "D36"; Value = "0.02533" },
    @{ Cell = "E36"; Value = "  +3.34%  " },
    @{ Cell = "D37"; Value = "5.485" },
    @{ Cell = "E37"; Value = "  +1.51%  " },
    @{ Cell = "D38"; Value = "0.06673" },
    @{ Cell = "E38"; Value = "  +3.94%  " },
    @{ Cell = "D39"; Value = "12.45" },
    @{ Cell = "E39"; Value = "  +9.24%  " },
    @{ Cell = "D40"; Value = "9.191" },
    @{ Cell = "E40"; Value = "  +5.17%  " },
    @{ Cell = "D41"; Value = "0.2215" },
    @{ Cell = "E41"; Value = "  +2.65%  " },
    @{ Cell = "D42"; Value = "0.6652" },
    @{ Cell = "E42"; Value = "  +2.03%  " },
    @{ Cell = "D43"; Value = "1.240" },
    @{ Cell = "E43"; Value = "  +1.74%  " },
    @{ Cell = "E44"; Value = "  +0.26%  " },
    @{ Cell = "D45"; Value = "13.72" },
    @{ Cell = "E45"; Value = "  +2.24%  " },
    @{ Cell = "E46"; Value = "  +1.37%  " },
    @{ Cell = "D47"; Value = "2.197" },
    @{ Cell = "E47"; Value = "  -1.85%  " },
    @{ Cell = "D48"; Value = "3.631" },
    @{ Cell = "E48"; Value = "  -0.22%  " },
    @{ Cell = "E49"; Value = "  +3.94%  " },
    @{ Cell = "D50"; Value = "124.91" },
    @{ Cell = "E50"; Value = "  +2.32%  " },
    @{ Cell = "D51"; Value = "81.17" },
    @{ Cell = "E51"; Value = "  +2.72%  " }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.Style = "Normal"
}
